$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# RUNMANAGER sheet1: flip the "execute" flag for the Health questionnaire test to "yes"
$ws1.Range("C6").Value = "yes"

# DATA sheet2: flip "execute" flag for the Health questionnaire / chrome row to "yes"
$ws2.Range("B2").Value = "yes"

# Append a new row (row 8) duplicating the Health questionnaire test, but for "edge" browser
$ws2.Range("A8").Value = "completeHealthQuestionnareTest"
$ws2.Range("B8").Value = "yes"
$ws2.Range("C8").Value = "edge"
$ws2.Range("D8").Value = "'98.0"
$ws2.Range("E8").Value = "bomaseko1"
$ws2.Range("F8").Value = "'QEtCVG9remFuMjAyMQ=="
$ws2.Range("G8").Value = "'Health"

# Update selections: DATA sheet selection moves to B8 (new row), no longer the active tab
$ws2.Range("B8").Select()

# RUNMANAGER sheet selection moves to C2 and becomes the active tab/sheet
$ws1.Range("C2").Select()
$ws1.Activate()
